$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark existing tasks as "done" in column C
$ws.Range("C5").Value = "done"
$ws.Range("C11").Value = "done"
$ws.Range("C12").Value = "done"
$ws.Range("C14").Value = "done"

# Add new rows: monsters + first items (materials and heals) + skills
$ws.Range("A25").Value = "merges"
$ws.Range("B25").Value = "Fabio(directeur cohésion)"

$ws.Range("A26").Value = "Bonus: Qui Sont Ils"
$ws.Range("B26").Value = "Roméo"

# Move selection to the newly added row, scrolled so row 15 is at the top
$ws.Range("B26").Select()
